$d = $word.ActiveDocument

$replacements = @(
    @("2024-08-29 Thursday", "2024-08-30 Friday"),
    @("370×3=1110", "489×4=1956"),
    @("811×3=2433", "394×3=1182"),
    @("800×2=1600", "460×4=1840"),
    @("228×3=684", "848×8=6784"),
    @("717×8=5736", "395×4=1580"),
    @("430×3=1290", "289×8=2312"),
    @("822×9=7398", "639×6=3834"),
    @("346×5=1730", "208×9=1872"),
    @("935×5=4675", "959×2=1918"),
    @("275×2=550", "721×8=5768"),
    @("830×6=4980", "138×2=276"),
    @("292×8=2336", "659×5=3295"),
    @("248×8=1984", "644×8=5152"),
    @("104×5=520", "468×7=3276"),
    @("801×9=7209", "970×2=1940"),
    @("450×2=900", "126×3=378"),
    @("345×2=690", "270×9=2430"),
    @("510×4=2040", "564×4=2256"),
    @("300×5=1500", "326×2=652"),
    @("778×6=4668", "195×3=585"),
    @("624×8=4992", "246×2=492"),
    @("193×8=1544", "304×9=2736"),
    @("218×2=436", "829×8=6632"),
    @("964×3=2892", "414×3=1242"),
    @("671×4=2684", "623×4=2492")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
